$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Update the title text
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Incorporating GenAI into Client Workflow"

# Replace the textbox's three paragraphs with a single paragraph of content
$s.Shapes.Item(2).TextFrame.TextRange.Text = "{'Introduction': {'content': 'Brief overview of GenAI and its benefits for clients.'}, 'Understanding GenAI': {'content': 'Explanation of what GenAI is and how it works.'}, 'Benefits of GenAI': {'content': 'Highlight the advantages of incorporating GenAI into workflow.'}, 'Integration Process': {'content': 'Step-by-step guide on how clients can integrate GenAI into their existing workflow.'}, 'Best Practices': {'content': 'Recommendations on how to maximize the effectiveness of GenAI in workflow.'}, 'Case Studies': {'content': 'Real-life examples of companies successfully using GenAI in their workflow.'}, 'Q&A': {'content': 'Open the floor for questions and provide answers to common queries.'}, 'Conclusion': {'content': 'Summarize key points and encourage clients to start incorporating GenAI into their workflow.'}}"
